# Generate Report for Handoff
#
# Refreshes the localization-status report for the new handoff run:
#   - the source markdown file's generated GUID-based name changed from
#     33b49b86-df5e-4034-acf5-45b2df68c05f.md to
#     5c230f97-e134-400e-bf0d-a37072738d9a.md
#   - the xliff handoff artifact hash changed from
#     31d082e952541c4d33c48013506941095de5abfb to
#     a87a23a00120cf53dfa44d5c7bd54114798f5ad8
#   - the associated timestamps moved forward to the new handoff run times

$wb = $excel.ActiveWorkbook

$oldGuidName = "33b49b86-df5e-4034-acf5-45b2df68c05f"
$newGuidName = "5c230f97-e134-400e-bf0d-a37072738d9a"

$oldHash = "31d082e952541c4d33c48013506941095de5abfb"
$newHash = "a87a23a00120cf53dfa44d5c7bd54114798f5ad8"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuidName.md"
$wsOverview.Range("B2").Value = "e2e\$newGuidName.md"
$wsOverview.Range("G2").Value = "2016-08-17 18:56:35"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuidName.md"
$wsZhCn.Range("G2").Value = "$newGuidName.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-17 18:56:30"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuidName.md"
$wsDeDe.Range("G2").Value = "$newGuidName.$newHash.de-de.xlf"
